# Auto-generated Excel COM-interop script that applies the
# Sophia_Profits market-price refresh described in the commit diff.
# For every touched cell we just assign the literal numeric value
# that the scheduled price-refresh run produced.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("H40").Value = 1003366.7
$ws.Range("I40").Value = 5050.5
$ws.Range("K40").Value = 5050.5
$ws.Range("M40").Value = -4875.5
$ws.Range("H42").Value = 1306.6666
$ws.Range("I42").Value = 420
$ws.Range("J42").Value = 1750
$ws.Range("K42").Value = 1260
$ws.Range("L42").Value = 5250
$ws.Range("M42").Value = -1030
$ws.Range("N42").Value = -5710
$ws.Range("H64").Value = 8000
$ws.Range("I64").Value = 8000
$ws.Range("K64").Value = 8000
$ws.Range("M64").Value = -7752
$ws.Range("H67").Value = 8000
$ws.Range("I67").Value = 8000
$ws.Range("K67").Value = 8000
$ws.Range("M67").Value = -7142
$ws.Range("H70").Value = 1363.6364
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 1200
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 3600
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -4140
$ws.Range("H73").Value = 1363.6364
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 1200
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 3600
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -5472
$ws.Range("H98").Value = 6666.1665
$ws.Range("I98").Value = 832.3333
$ws.Range("J98").Value = 12500
$ws.Range("K98").Value = 832.3333
$ws.Range("L98").Value = 12500
$ws.Range("M98").Value = 665.6667
$ws.Range("N98").Value = -15496
$ws.Range("H103").Value = 1176.3846
$ws.Range("I103").Value = 323.5
$ws.Range("J103").Value = 1555.4445
$ws.Range("K103").Value = 970.5
$ws.Range("L103").Value = 4666.333500000001
$ws.Range("M103").Value = -384.5
$ws.Range("N103").Value = -5838.333500000001
$ws.Range("H111").Value = 4442
$ws.Range("I111").Value = 1794
$ws.Range("J111").Value = 5766
$ws.Range("K111").Value = 5382
$ws.Range("L111").Value = 17298
$ws.Range("M111").Value = -2315
$ws.Range("N111").Value = -23432
$ws.Range("H122").Value = 6666.1665
$ws.Range("I122").Value = 832.3333
$ws.Range("J122").Value = 12500
$ws.Range("K122").Value = 2496.9999
$ws.Range("L122").Value = 37500
$ws.Range("M122").Value = -46.9998999999998
$ws.Range("N122").Value = -42400
$ws.Range("H138").Value = 4060.077
$ws.Range("J138").Value = 5081.722
$ws.Range("L138").Value = 15245.166
$ws.Range("N138").Value = -25525.166
$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 1500
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4500
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = 680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1600
$ws.Range("J2").Value = 768
$ws.Range("L2").Value = 768
$ws.Range("N2").Value = -994
$ws.Range("H32").Value = 4963.5
$ws.Range("I32").Value = 4963.5
$ws.Range("K32").Value = 4963.5
$ws.Range("M32").Value = -4676.5
$ws.Range("H74").Value = 9954.323
$ws.Range("I74").Value = 11076.444
$ws.Range("K74").Value = 11076.444
$ws.Range("M74").Value = -10202.444
$ws.Range("H76").Value = 59999
$ws.Range("J76").Value = 59999
$ws.Range("L76").Value = 59999
$ws.Range("N76").Value = -60675
$ws.Range("H77").Value = 9954.323
$ws.Range("I77").Value = 11076.444
$ws.Range("K77").Value = 55382.22
$ws.Range("M77").Value = -51014.22
$ws.Range("H79").Value = 59999
$ws.Range("J79").Value = 59999
$ws.Range("L79").Value = 59999
$ws.Range("N79").Value = -62339
$ws.Range("H109").Value = 115188.5
$ws.Range("J109").Value = 115188.5
$ws.Range("L109").Value = 115188.5
$ws.Range("N109").Value = -117962.5
$ws.Range("H116").Value = 1600
$ws.Range("J116").Value = 768
$ws.Range("L116").Value = 768
$ws.Range("N116").Value = -5356

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1600
$ws.Range("J3").Value = 768
$ws.Range("L3").Value = 768
$ws.Range("N3").Value = -996
$ws.Range("H20").Value = 4467.75
$ws.Range("I20").Value = 1088.4286
$ws.Range("J20").Value = 9198.799999999999
$ws.Range("K20").Value = 1088.4286
$ws.Range("L20").Value = 9198.799999999999
$ws.Range("M20").Value = -841.4286
$ws.Range("N20").Value = -9692.799999999999
$ws.Range("H99").Value = 2084.8333
$ws.Range("I99").Value = 1601.8
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 1601.8
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -103.8
$ws.Range("N99").Value = -7496
$ws.Range("H107").Value = 665.3333
$ws.Range("I107").Value = 648
$ws.Range("K107").Value = 648
$ws.Range("M107").Value = 1272

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4611.769
$ws.Range("I58").Value = 4579.4165
$ws.Range("K58").Value = 4579.4165
$ws.Range("M58").Value = -4376.4165
$ws.Range("H59").Value = 60130
$ws.Range("J59").Value = 60130
$ws.Range("L59").Value = 60130
$ws.Range("N59").Value = -62420
$ws.Range("H86").Value = 7996.143
$ws.Range("J86").Value = 7748.5
$ws.Range("L86").Value = 7748.5
$ws.Range("N86").Value = -9994.5
$ws.Range("H89").Value = 7996.143
$ws.Range("J89").Value = 7748.5
$ws.Range("L89").Value = 38742.5
$ws.Range("N89").Value = -49974.5
$ws.Range("H105").Value = 3428.5715
$ws.Range("I105").Value = 3192.3076
$ws.Range("J105").Value = 3812.5
$ws.Range("K105").Value = 3192.3076
$ws.Range("L105").Value = 3812.5
$ws.Range("M105").Value = -1445.3076
$ws.Range("N105").Value = -7306.5
$ws.Range("H134").Value = 6761.0713
$ws.Range("I134").Value = 6667.6
$ws.Range("J134").Value = 6994.75
$ws.Range("K134").Value = 20002.8
$ws.Range("L134").Value = 20984.25
$ws.Range("M134").Value = -17467.8
$ws.Range("N134").Value = -26054.25
$ws.Range("H136").Value = 4611.769
$ws.Range("I136").Value = 4579.4165
$ws.Range("K136").Value = 13738.2495
$ws.Range("M136").Value = -11188.2495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2749
$ws.Range("J131").Value = 2687.5
$ws.Range("L131").Value = 8062.5
$ws.Range("N131").Value = -18142.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 253
$ws.Range("I2").Value = 207.25
$ws.Range("K2").Value = 207.25
$ws.Range("M2").Value = -94.25
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10346
$ws.Range("H80").Value = 6598.5713
$ws.Range("I80").Value = 2737.5
$ws.Range("J80").Value = 11746.667
$ws.Range("K80").Value = 2737.5
$ws.Range("L80").Value = 11746.667
$ws.Range("M80").Value = -1739.5
$ws.Range("N80").Value = -13742.667
$ws.Range("H83").Value = 6598.5713
$ws.Range("I83").Value = 2737.5
$ws.Range("J83").Value = 11746.667
$ws.Range("K83").Value = 13687.5
$ws.Range("L83").Value = 58733.335
$ws.Range("M83").Value = -8695.5
$ws.Range("N83").Value = -68717.33499999999
$ws.Range("H132").Value = 4364.2
$ws.Range("I132").Value = 3941
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 11823
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -9293
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3170.9092
$ws.Range("I68").Value = 3248
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 3248
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = -2499
$ws.Range("N68").Value = -3898
$ws.Range("H71").Value = 3170.9092
$ws.Range("I71").Value = 3248
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 16240
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = -12496
$ws.Range("N71").Value = -19488
$ws.Range("H132").Value = 3051.1333
$ws.Range("I132").Value = 2251.818
$ws.Range("K132").Value = 6755.454000000001
$ws.Range("M132").Value = -4225.454000000001
$ws.Range("H136").Value = 3498.5
$ws.Range("I136").Value = 3498.5
$ws.Range("K136").Value = 10495.5
$ws.Range("M136").Value = -7945.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1732.3334
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H136").Value = 3477.9333
$ws.Range("I136").Value = 3369.2144
$ws.Range("K136").Value = 10107.6432
$ws.Range("M136").Value = -7557.643199999999

Write-Output "Sophia_Profits price refresh applied"
